$d = $word.ActiveDocument

# 1) LOQ4031 -> LOQ4010 (text replacement within the existing run)
$d.Content.Find.Execute(
    "LOQ4031 -  Química Geral I  (Requisito)", $true, $false, $false, $false, $false,
    $true, 1, $false, "LOQ4010 -  Introdução à  Engenharia  Química  (Requisito)", 2
) | Out-Null

# 2) Remove the LOQ4073 line entirely (run text + its trailing line break)
$rng = $d.Content
$rng.Find.Execute("LOQ4073 -  Química Geral II  (Requisito)") | Out-Null
$toDelete = $d.Range($rng.Start, $rng.End + 1)
$toDelete.Delete()

# 3) Insert two new requisito lines right after the LOQ4095 line
$rng2 = $d.Content
$rng2.Find.Execute("LOQ4095 -  Química Geral Experimental  (Requisito)") | Out-Null
$insPoint = $d.Range($rng2.End + 1, $rng2.End + 1)
$insPoint.InsertAfter("LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito)" + [char]11)
$insPoint2 = $d.Range($insPoint.End, $insPoint.End)
$insPoint2.InsertAfter("LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)" + [char]11)
